$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New volunteer entry (row 50). Write the email (column B) before the name
# (column A) so new shared-string entries land in the same order as the
# source edit: "aritankovic@gmail.com" then "Armin Tankovic".
$ws.Range("B50").Value = "aritankovic@gmail.com"
$ws.Range("A50").Value = "Armin Tankovic"
$ws.Range("D50").Value = "UTM"
$ws.Range("E50").Value = 6478227275

# Scroll the sheet down and move the active selection, matching the
# author's final view state.
$ws.Range("E54").Select() | Out-Null
